$wb = $excel.ActiveWorkbook

$newNames = @("summ33331918", "summ33571800", "summ33889263", "summ34210050", "summ34506118", "summ34854272", "summ35164390", "summ35478906", "summ35827530", "summ36142245", "summ36435647", "summ36748617", "summ37105849", "summ37451580", "summ37777217", "summ38090097", "summ38438597", "summ38770699", "summ39078610", "summ39408486", "summ39707177", "summ40020020", "summ40314839", "summ40622380", "summ40910309", "summ41233320", "summ41553723", "summ41887608", "summ42178743", "summ42491520", "summ42791469", "summ43189484", "summ43528482", "summ43811131", "summ44197558", "summ44757364", "summ45122734", "summ45591436", "summ45898905", "summ46189099", "summ46492096", "summ46773965", "summ47065431", "summ47356762", "summ47655533", "summ47960775", "summ48269053", "summ48564749", "summ48871006", "summ49177697")

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
